# Append EcoTech coupon rows (2-9) to the active sheet, matching the
# "company, code, points, user_email, created_at" header already in row 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Company = "EcoTech"; Code = "EcoTech-CZJ734702"; Points = 100; Created = 45752.21452033565 },
    @{ Company = "EcoTech"; Code = "EcoTech-8RW737589"; Points = 100; Created = 45752.21455285879 },
    @{ Company = "EcoTech"; Code = "EcoTech-4RW738741"; Points = 100; Created = 45752.21456609954 },
    @{ Company = "EcoTech"; Code = "EcoTech-OQ5739494"; Points = 100; Created = 45752.21457478009 },
    @{ Company = "EcoTech"; Code = "EcoTech-GZR739976"; Points = 100; Created = 45752.21458042824 },
    @{ Company = "EcoTech"; Code = "EcoTech-FYP740533"; Points = 100; Created = 45752.2145868287 },
    @{ Company = "EcoTech"; Code = "EcoTech-CQ6740971"; Points = 100; Created = 45752.21459194444 },
    @{ Company = "EcoTech"; Code = "EcoTech-8M3741255"; Points = 100; Created = 45752.21459517183 }
)

$firstRow = 2
$lastRow = $firstRow + $rows.Count - 1

$r = $firstRow
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.Company
    $ws.Cells.Item($r, 2).Value = $row.Code
    $ws.Cells.Item($r, 3).Value = $row.Points
    $ws.Cells.Item($r, 5).Value = $row.Created
    $r = $r + 1
}

# Give the first created_at cell its timestamp number format. Registering
# the lowercase variant first (numFmtId 164) and then overwriting the same
# cell's format with the uppercase variant (numFmtId 165) mirrors how the
# source file was produced: both format codes end up in the stylesheet,
# while the cell itself lands on the 165 style.
$firstDateCell = $ws.Cells.Item($firstRow, 5)
$firstDateCell.NumberFormat = "yyyy-mm-dd h:mm:ss"
$firstDateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Propagate that exact style (not a freshly-resolved one) to the rest of the
# created_at column via copy/paste-special so every row shares a single
# stylesheet entry instead of each cell minting its own.
if ($lastRow -gt $firstRow) {
    $firstDateCell.Copy()
    $ws.Range($ws.Cells.Item($firstRow + 1, 5), $ws.Cells.Item($lastRow, 5)).PasteSpecial(-4122)
}
